$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C144").Value = "未完成"
$ws.Range("C145").Value = "未完成"
$ws.Range("C146").Value = "用户和宝贝信息添加完成"
$ws.Range("C147").Value = "已完成"
$ws.Range("C148").Value = "未完成"
$ws.Range("C149").Value = "求购信息添加完成"

$ws.Range("A150").Value = "总结：应该再细分，提高工作效率"

$ws.Rows("144").RowHeight = 22.5
$ws.Rows("147").RowHeight = 22.5
$ws.Rows("148").RowHeight = 22.5

$ws.Range("A150").Select()
